# Fruta / hortaliza, semanal
# Insert a new weekly price-report row above the current row 93 on the
# active sheet (Comercializadora del Agro de Limarí - Haba), shifting all
# subsequent rows down by one, and populate the new row with this week's
# data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 93 (and everything below it) down one row.
$ws.Rows.Item(93).Insert()

# Fill in the new row 93 with the latest weekly record.
$ws.Cells.Item(93, 1).Value = 2
$ws.Cells.Item(93, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(93, 3).Value = "Coquimbo"
$ws.Cells.Item(93, 4).Value = 45204
$ws.Cells.Item(93, 5).Value = 4
$ws.Cells.Item(93, 6).Value = 100112026
$ws.Cells.Item(93, 7).Value = "Haba"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 1100
$ws.Cells.Item(93, 11).Value = 7000
$ws.Cells.Item(93, 12).Value = 8000
$ws.Cells.Item(93, 13).Value = 7500
$ws.Cells.Item(93, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(93, 16).Value = 300
$ws.Cells.Item(93, 17).Value = 25
$ws.Cells.Item(93, 18).Value = "Hortaliza"
